# Dash improved, sword slash init
# Adds two new "idea" rows to the tracker sheet: "Nero inspired sword power"
# and "Shooting from the sword".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "Nero inspired sword power" ------------------------------------
# Same look as the existing rows 2/3: full thin-box border, centered;
# A/B wrap text, C does not. Copy the formatting straight from row 2 so the
# new row reuses the very same style entries instead of creating new ones.
$ws.Range("A2:B2").Copy()
$ws.Range("A4:B4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("C4").PasteSpecial(-4122)     # xlPasteFormats

$ws.Range("A4").Value = "Nero inspired sword power"
$ws.Range("B4").Value = "As Travis holds a button, he will charge his sword with the equipped power (make him twirl the sword once he has the attacked charged up), making it light on fire, be covered in electricity, ice, etc. This effect will also be achieved by pressing the same button just after an attack lands"
$ws.Range("C4").Value = "Not started"

$ws.Rows.Item(4).RowHeight = 120

# --- Row 5: "Shooting from the sword" --------------------------------------
# Same centered look, but only left/right thin borders (no top/bottom).
# Start from row 2's formatting (full box) and strip the top/bottom edges on
# each cell.
$ws.Range("A2:B2").Copy()
$ws.Range("A5:B5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)     # xlPasteFormats

foreach ($addr in @("A5", "B5", "C5")) {
    $cell = $ws.Range($addr)
    $cell.Borders.Item(8).LineStyle = -4142  # xlEdgeTop -> xlLineStyleNone
    $cell.Borders.Item(9).LineStyle = -4142  # xlEdgeBottom -> xlLineStyleNone
}

$ws.Range("A5").Value = "Shooting from the sword"
$ws.Range("B5").Value = "When you shoot any power, it will come out of the sword, instead of the hand (if not fighting with the sword), Travis will twirl his sword to achieve this."
$ws.Range("C5").Value = "Not Started"

$ws.Rows.Item(5).RowHeight = 60

# --- Sheet view: land the selection on B5 -----------------------------------
$ws.Range("B5").Select() | Out-Null

Write-Host "Added rows 4-5 (Nero inspired sword power / Shooting from the sword)"
